$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 528, shifting the existing rows 528:536 down to 531:539.
$ws.Rows("528:530").Insert()

# Row 528 (new): Lechuga, Conconina(o), Primera - Fecha 03-02-2022
$ws.Cells.Item(528, 1).Value = 7
$ws.Cells.Item(528, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(528, 3).Value = "Ñuble"
$ws.Cells.Item(528, 4).Value = 44595
$ws.Cells.Item(528, 5).Value = 16
$ws.Cells.Item(528, 6).Value = 100112033
$ws.Cells.Item(528, 7).Value = "Lechuga"
$ws.Cells.Item(528, 8).Value = "Conconina(o)"
$ws.Cells.Item(528, 9).Value = "Primera"
$ws.Cells.Item(528, 10).Value = 120
$ws.Cells.Item(528, 11).Value = 4500
$ws.Cells.Item(528, 12).Value = 5000
$ws.Cells.Item(528, 13).Value = 4750
$ws.Cells.Item(528, 14).Value = "$/caja 10 unidades"
$ws.Cells.Item(528, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(528, 16).Value = 475
$ws.Cells.Item(528, 17).Value = 10
$ws.Cells.Item(528, 18).Value = "Hortaliza"

# Row 529 (new): Lechuga, Escarola, Primera - Fecha 03-02-2022
$ws.Cells.Item(529, 1).Value = 7
$ws.Cells.Item(529, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(529, 3).Value = "Ñuble"
$ws.Cells.Item(529, 4).Value = 44595
$ws.Cells.Item(529, 5).Value = 16
$ws.Cells.Item(529, 6).Value = 100112033
$ws.Cells.Item(529, 7).Value = "Lechuga"
$ws.Cells.Item(529, 8).Value = "Escarola"
$ws.Cells.Item(529, 9).Value = "Primera"
$ws.Cells.Item(529, 10).Value = 120
$ws.Cells.Item(529, 11).Value = 5000
$ws.Cells.Item(529, 12).Value = 5500
$ws.Cells.Item(529, 13).Value = 5250
$ws.Cells.Item(529, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(529, 15).Value = "Región del Maule"
$ws.Cells.Item(529, 16).Value = 350
$ws.Cells.Item(529, 17).Value = 15
$ws.Cells.Item(529, 18).Value = "Hortaliza"

# Row 530 (new): Lechuga, Marina, Primera - Fecha 03-02-2022
$ws.Cells.Item(530, 1).Value = 7
$ws.Cells.Item(530, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(530, 3).Value = "Ñuble"
$ws.Cells.Item(530, 4).Value = 44595
$ws.Cells.Item(530, 5).Value = 16
$ws.Cells.Item(530, 6).Value = 100112033
$ws.Cells.Item(530, 7).Value = "Lechuga"
$ws.Cells.Item(530, 8).Value = "Marina"
$ws.Cells.Item(530, 9).Value = "Primera"
$ws.Cells.Item(530, 10).Value = 120
$ws.Cells.Item(530, 11).Value = 5000
$ws.Cells.Item(530, 12).Value = 5500
$ws.Cells.Item(530, 13).Value = 5250
$ws.Cells.Item(530, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(530, 15).Value = "Región del Maule"
$ws.Cells.Item(530, 16).Value = 292
$ws.Cells.Item(530, 17).Value = 18
$ws.Cells.Item(530, 18).Value = "Hortaliza"
